$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current values of rows 2-23, columns A-F, before overwriting
# anything (the new layout is a permutation of the existing rows).
$snapshot = @{}
for ($r = 2; $r -le 23; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le 6; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping of new row number -> source row number (from the old snapshot)
$mapping = @{
    2  = 10
    3  = 14
    4  = 15
    5  = 8
    6  = 13
    7  = 2
    8  = 7
    9  = 3
    10 = 6
    11 = 11
    12 = 12
    13 = 4
    14 = 9
    15 = 5
    16 = 20
    17 = 16
    18 = 18
    19 = 21
    20 = 19
    21 = 17
    22 = 23
    23 = 22
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $srcVals[$c]
    }
}
